$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying edit is a cyclic re-shuffle of the per-observation species
# data (and, for a few rows, the date/"Bestämningsmetod" metadata) among the
# data rows. Row numbers stay the same; only the species-related payload
# moves between rows, following these cycles:
#   2 <- 5 <- 3 <- 2
#   6 <- 9 <- 6
#   8 <- 11 <- 10 <- 8

# --- Row 2 (was row 5's species data) ---
$ws.Range("A2").Value = 111780621
$ws.Range("B2").Value = 56543
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 103021
$ws.Range("F2").Value = "Talltita"
$ws.Range("G2").Value = "Poecile montanus"
$ws.Range("H2").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q2").Value = 707631.1509720345
$ws.Range("R2").Value = 7397277.54798521

# --- Row 3 (was row 2's species data) ---
$ws.Range("A3").Value = 111780624
$ws.Range("B3").Value = 95532
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 221945
$ws.Range("F3").Value = "Revlummer"
$ws.Range("G3").Value = "Lycopodium annotinum"
$ws.Range("H3").Value = "L."
$ws.Range("Q3").Value = 707600.9335272597
$ws.Range("R3").Value = 7397313.141869167

# --- Row 5 (was row 3's species data) ---
$ws.Range("A5").Value = 111780627
$ws.Range("B5").Value = 78604
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 6461
$ws.Range("F5").Value = "Norrlandslav"
$ws.Range("G5").Value = "Nephroma arcticum"
$ws.Range("H5").Value = "(L.) Torss."
$ws.Range("Q5").Value = 707647.2196405758
$ws.Range("R5").Value = 7397286.731778639

# --- Row 6 (was row 9's species data) ---
$ws.Range("A6").Value = 111816132
$ws.Range("B6").Value = 95532
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221945
$ws.Range("F6").Value = "Revlummer"
$ws.Range("G6").Value = "Lycopodium annotinum"
$ws.Range("H6").Value = "L."
$ws.Range("Q6").Value = 707589.6730983062
$ws.Range("R6").Value = 7397240.139162621

# --- Row 8 (was row 11's species data + row 11's date metadata) ---
$ws.Range("A8").Value = 111816119
$ws.Range("B8").Value = 56543
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 103021
$ws.Range("F8").Value = "Talltita"
$ws.Range("G8").Value = "Poecile montanus"
$ws.Range("H8").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q8").Value = 707595.5401507822
$ws.Range("R8").Value = 7397262.905378895
$ws.Range("Y8").Value = "2023-08-22"
$ws.Range("AA8").Value = "2023-08-22"
$ws.Range("AF8").ClearContents()

# --- Row 9 (was row 6's species data) ---
$ws.Range("A9").Value = 111816118
$ws.Range("B9").Value = 78107
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6453
$ws.Range("F9").Value = "Vedskivlav"
$ws.Range("G9").Value = "Hertelidea botryosa"
$ws.Range("H9").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q9").Value = 707670.4513803272
$ws.Range("R9").Value = 7397327.948038339

# --- Row 10 (was row 8's species data + row 8's date metadata) ---
$ws.Range("A10").Value = 111816142
$ws.Range("B10").Value = 78604
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 6461
$ws.Range("F10").Value = "Norrlandslav"
$ws.Range("G10").Value = "Nephroma arcticum"
$ws.Range("H10").Value = "(L.) Torss."
$ws.Range("Q10").Value = 707613.3456041727
$ws.Range("R10").Value = 7397270.22663033
$ws.Range("Y10").Value = "2023-08-29"
$ws.Range("AA10").Value = "2023-08-29"

# --- Row 11 (was row 10's species data + row 10's empty "Bestämningsmetod" cell) ---
$ws.Range("A11").Value = 111816145
$ws.Range("B11").Value = 77597
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 864
$ws.Range("F11").Value = "Knottrig blåslav"
$ws.Range("G11").Value = "Hypogymnia bitteri"
$ws.Range("H11").Value = "(Lynge) Ahti"
$ws.Range("Q11").Value = 707626.9948496711
$ws.Range("R11").Value = 7397311.517900761
$ws.Range("AF10").Copy($ws.Range("AF11"))
